$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple text/text-like updates (safe from numeric auto-conversion) ---
$ws.Range("D2").Value = "64.569.59"
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("D3").Value = "3.153.85"
$ws.Range("E3").Value = "  -0.91%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("E5").Value = "  +1.76%  "
$ws.Range("E6").Value = "  -2.86%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").Value = "3.149.81"
$ws.Range("E8").Value = "  -0.97%  "
$ws.Range("E9").Value = "  -1.26%  "
$ws.Range("E10").Value = "  -1.79%  "
$ws.Range("E11").Value = "  -2.52%  "
$ws.Range("E12").Value = "  -1.21%  "
$ws.Range("E13").Value = "  -1.29%  "
$ws.Range("E14").Value = "  -4.00%  "
$ws.Range("D15").Value = "3.681.06"
$ws.Range("E15").Value = "  -0.57%  "
$ws.Range("E16").Value = "  +2.70%  "
$ws.Range("D17").Value = "64.562.20"
$ws.Range("E17").Value = "  -0.52%  "
$ws.Range("D18").Value = "3.159.98"
$ws.Range("E18").Value = "  -0.95%  "
$ws.Range("E19").Value = "  -2.26%  "
$ws.Range("E20").Value = "  -0.96%  "
$ws.Range("E21").Value = "  -1.09%  "
$ws.Range("B22").Value = "Polygon"
$ws.Range("C22").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("E23").Value = "  +2.75%  "
$ws.Range("E24").Value = "  -2.58%  "
$ws.Range("E25").Value = "  -1.18%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  -3.39%  "
$ws.Range("E28").Value = "  -2.36%  "
$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("E29").Value = "  -0.20%  "
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("E30").Value = "  -1.24%  "
$ws.Range("E31").Value = "  -7.54%  "
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("E32").Value = "  +0.21%  "
$ws.Range("B33").Value = "Stacks"
$ws.Range("C33").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("E33").Value = "  -0.26%  "
$ws.Range("E34").Value = "  -2.33%  "
$ws.Range("E35").Value = "  +1.46%  "
$ws.Range("D36").Value = "0.0₃0789"
$ws.Range("E36").Value = "  +6.54%  "
$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("E37").Value = "  -2.52%  "
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("E38").Value = "  +0.14%  "
$ws.Range("E39").Value = "  -2.73%  "
$ws.Range("E40").Value = "  -0.94%  "
$ws.Range("E41").Value = "  -1.18%  "
$ws.Range("E42").Value = "  -5.90%  "
$ws.Range("E43").Value = "  -1.82%  "
$ws.Range("D44").Value = "2.849.84"
$ws.Range("E44").Value = "  -2.65%  "
$ws.Range("E45").Value = "  -5.39%  "
$ws.Range("E46").Value = "  -3.27%  "
$ws.Range("E47").Value = "  +3.74%  "
$ws.Range("E48").Value = "  -3.23%  "
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("E50").Value = "  -2.10%  "
$ws.Range("E51").Value = "  -0.19%  "

# --- Numeric-looking text values: force text format to avoid Excel converting them to numbers ---
# Preserve original style, apply text format, set value, restore style to avoid leaving stray formatting diffs
$cellsToFix = @(
    @{Cell="D5"; Value="613.38"},
    @{Cell="D6"; Value="148.22"},
    @{Cell="D9"; Value="0.527"},
    @{Cell="D10"; Value="0.152"},
    @{Cell="D11"; Value="5.45"},
    @{Cell="D12"; Value="0.472"},
    @{Cell="D13"; Value="0.0000259"},
    @{Cell="D14"; Value="35.66"},
    @{Cell="D19"; Value="6.92"},
    @{Cell="D20"; Value="480.70"},
    @{Cell="D21"; Value="14.70"},
    @{Cell="D22"; Value="0.719"},
    @{Cell="D23"; Value="8.03"},
    @{Cell="D24"; Value="13.69"},
    @{Cell="D25"; Value="83.93"},
    @{Cell="D27"; Value="2.84"},
    @{Cell="D28"; Value="8.57"},
    @{Cell="D29"; Value="0.121"},
    @{Cell="D30"; Value="7.05"},
    @{Cell="D31"; Value="2.10"},
    @{Cell="D32"; Value="1.00"},
    @{Cell="D33"; Value="2.72"},
    @{Cell="D34"; Value="26.40"},
    @{Cell="D35"; Value="1.14"},
    @{Cell="D37"; Value="6.02"},
    @{Cell="D38"; Value="3.24"},
    @{Cell="D39"; Value="53.15"},
    @{Cell="D40"; Value="459.41"},
    @{Cell="D41"; Value="0.0402"},
    @{Cell="D42"; Value="0.120"},
    @{Cell="D43"; Value="8.39"},
    @{Cell="D45"; Value="2.34"},
    @{Cell="D46"; Value="0.269"},
    @{Cell="D47"; Value="2.45"},
    @{Cell="D48"; Value="26.57"},
    @{Cell="D49"; Value="1.00"},
    @{Cell="D50"; Value="0.114"},
    @{Cell="D51"; Value="120.30"}
)

foreach ($item in $cellsToFix) {
    $rng = $ws.Range($item.Cell)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $item.Value
    $rng.Style = $origStyle
}
